$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-16 shift up to rows 12-15 (Trae Young, Terry Rozier, Dalton Knecht, Donte DiVincenzo)
$ws.Range("A12").Value = "Trae Young"
$ws.Range("B12").Value = "PG"
$ws.Range("C12").Value = "Atlanta Hawks"

$ws.Range("A13").Value = "Terry Rozier"
$ws.Range("B13").Value = "PG"
$ws.Range("C13").Value = "Miami Heat"

$ws.Range("A14").Value = "Dalton Knecht"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Los Angeles Lakers"

$ws.Range("A15").Value = "Donte DiVincenzo"
$ws.Range("B15").Value = "SG,SF"
$ws.Range("C15").Value = "Minnesota Timberwolves"

# Row 16 becomes what was row 19 (Nicolas Claxton, C, Brooklyn Nets)
$ws.Range("A16").Value = "Nicolas Claxton"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Brooklyn Nets"

# Row 19 becomes what was row 12 (Norman Powell, SG,SF, LA Clippers)
$ws.Range("A19").Value = "Norman Powell"
$ws.Range("B19").Value = "SG,SF"
$ws.Range("C19").Value = "LA Clippers"
